$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.03%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.097"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08132"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.18%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.967"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.75%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.933"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.99%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9281"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.71%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1443"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "12.59%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1952"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.91%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09099"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.70%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03504"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.80%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09831"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.15%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001404"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006147"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.57%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.630"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.37%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.202"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.31%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.64%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.74%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.824"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.89%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.35%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04434"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.10%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.32%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004847"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "5.12%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.02%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02100"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "7.79%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05122"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.42%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007488"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01015"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.28%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.97%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.36%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009428"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.13%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006228"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.48%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.01%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003060"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001602"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.01%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.01%"
